$wb = $excel.ActiveWorkbook

# --- Functions sheet (sheet2): add the new Day-3 "Zig Zag" example row ---
$wsFunctions = $wb.Worksheets.Item("Functions")

$wsFunctions.Range("A18").Value = "print(""Welcome to"" , end = ' ') " + [char]10 + "print(""GeeksforGeeks"", end = ' ')"
$wsFunctions.Range("A18").WrapText = $true
$wsFunctions.Range("B18").Value = "# ends the output with a <space> "
$wsFunctions.Range("C18").Value = "Welcome to GeeksforGeeks"

$wsFunctions.Rows.Item(18).RowHeight = 28.8

# Widen column A slightly to fit the new content
$wsFunctions.Columns.Item(1).ColumnWidth = 29.65

# --- Update the active sheet / selection state ---
# Previously "Exceptions" (sheet 4) was the active/selected sheet with B3
# selected; now it moves to "Functions" (sheet 2), and Exceptions' own
# selection resets to B1.
$wsExceptions = $wb.Worksheets.Item("Exceptions")
$wsExceptions.Activate()
$wsExceptions.Range("B1").Select()

$wsFunctions.Activate()
$wsFunctions.Range("B21").Select()
